$d = $word.ActiveDocument

$pairs = @(
    ,@("81+1=82", "82-77=5")
    ,@("29-19=10", "33+35=68")
    ,@("88-82=6", "98-45=53")
    ,@("52-8=44", "35-0=35")
    ,@("88+0=88", "53+6=59")
    ,@("89-82=7", "11+59=70")
    ,@("26+24=50", "21+21=42")
    ,@("40-5=35", "58-13=45")
    ,@("7+46=53", "31+12=43")
    ,@("44+29=73", "32+40=72")
    ,@("82-9=73", "5+1=6")
    ,@("79+4=83", "39+2=41")
    ,@("87+10=97", "62-28=34")
    ,@("34+42=76", "19+34=53")
    ,@("85-64=21", "78-31=47")
    ,@("59-47=12", "16+80=96")
    ,@("38+56=94", "60-54=6")
    ,@("12+58=70", "8+82=90")
    ,@("24+18=42", "97+1=98")
    ,@("63+35=98", "54-24=30")
    ,@("88-66=22", "8+66=74")
    ,@("40-26=14", "71+10=81")
    ,@("27+24=51", "9+60=69")
    ,@("49+45=94", "4+12=16")
    ,@("52-47=5", "44-23=21")
    ,@("84+10=94", "81+6=87")
    ,@("38-13=25", "32-2=30")
    ,@("41+38=79", "70-14=56")
    ,@("22+50=72", "99-35=64")
    ,@("74-44=30", "14-0=14")
    ,@("70-61=9", "21+14=35")
    ,@("42+14=56", "73-44=29")
    ,@("76+0=76", "69-53=16")
    ,@("62+27=89", "63+29=92")
    ,@("89-84=5", "33+47=80")
    ,@("65-46=19", "98-27=71")
    ,@("20-0=20", "12-10=2")
    ,@("11+16=27", "52+27=79")
    ,@("68+26=94", "78-38=40")
    ,@("65+6=71", "68+6=74")
    ,@("10+45=55", "45-11=34")
    ,@("10+23=33", "46+22=68")
    ,@("82-27=55", "51-42=9")
    ,@("69-18=51", "11+4=15")
    ,@("40-35=5", "24+74=98")
    ,@("32+16=48", "2+10=12")
    ,@("3+80=83", "30+21=51")
    ,@("43+4=47", "52+33=85")
    ,@("27+71=98", "61+16=77")
    ,@("44-26=18", "61+13=74")
    ,@("77-70=7", "1+76=77")
    ,@("24-9=15", "10+1=11")
    ,@("83-36=47", "38+41=79")
    ,@("78-10=68", "39+3=42")
    ,@("67+25=92", "72+26=98")
    ,@("32-10=22", "1+27=28")
    ,@("35+7=42", "46+16=62")
    ,@("18+19=37", "18-12=6")
    ,@("7+24=31", "7+16=23")
    ,@("67-43=24", "11+81=92")
    ,@("3-3=0", "14+4=18")
    ,@("3+17=20", "57+34=91")
    ,@("72-15=57", "26+43=69")
    ,@("34+12=46", "32-31=1")
    ,@("1+58=59", "72-0=72")
    ,@("19+57=76", "60+20=80")
    ,@("48+33=81", "97-28=69")
    ,@("25+29=54", "29+13=42")
    ,@("94-49=45", "5+72=77")
    ,@("32+63=95", "52-0=52")
    ,@("38+48=86", "4+5=9")
    ,@("24+6=30", "8+59=67")
    ,@("51+15=66", "25+13=38")
    ,@("32+49=81", "87-1=86")
    ,@("2+18=20", "81-17=64")
    ,@("55-14=41", "94-15=79")
    ,@("44-6=38", "42+32=74")
    ,@("46+5=51", "54+17=71")
    ,@("0+48=48", "95-51=44")
    ,@("83+6=89", "74-19=55")
    ,@("67-23=44", "48-4=44")
    ,@("28+35=63", "45+47=92")
    ,@("23+61=84", "89-51=38")
    ,@("9+49=58", "73-37=36")
    ,@("80-55=25", "27+12=39")
    ,@("73+17=90", "4+29=33")
    ,@("14+9=23", "51-2=49")
    ,@("55+0=55", "24+51=75")
    ,@("49-19=30", "9+75=84")
    ,@("99-67=32", "42+19=61")
    ,@("19+2=21", "41+32=73")
    ,@("23-1=22", "36+39=75")
    ,@("63-54=9", "86-11=75")
    ,@("28+3=31", "37+0=37")
    ,@("7+57=64", "64+27=91")
    ,@("53+10=63", "14-4=10")
    ,@("99-53=46", "55+29=84")
    ,@("38+29=67", "41-37=4")
    ,@("52-44=8", "61-14=47")
    ,@("73-62=11", "78-6=72")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $null = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done: replaced $($pairs.Count) cells"
